$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 94, shifting the existing rows 94..158 down to 95..159.
$ws.Rows("94:94").Insert()

# Populate the newly inserted row 94 with this week's new record.
$ws.Range("A94").Value = 10
$ws.Range("B94").Value = 'Vega Modelo de Temuco'
$ws.Range("C94").Value = 'La Araucanía'
$ws.Range("D94").Value = 44827
$ws.Range("E94").Value = 9
$ws.Range("F94").Value = 100112031
$ws.Range("G94").Value = 'Poroto verde'
$ws.Range("H94").Value = 'Sin especificar'
$ws.Range("I94").Value = 'Primera'
$ws.Range("J94").Value = 30
$ws.Range("K94").Value = 32000
$ws.Range("L94").Value = 32000
$ws.Range("M94").Value = 32000
$ws.Range("N94").Value = '$/malla 25 kilos'
$ws.Range("O94").Value = 'Provincia de Limarí'
$ws.Range("P94").Value = 1280
$ws.Range("Q94").Value = 25
$ws.Range("R94").Value = 'Hortaliza'
